$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

# 1. Rename header row column names: "_old" -> "_FV2404", "_new" -> "_FV2410"
#    (headers were suffixed with the generic "_old"/"_new" markers; they now
#    need to carry the concrete format-version names they were diffed from/to)
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value2
    if ($text -ne $null) {
        if ($text.EndsWith("_old")) {
            $cell.Value2 = $text.Substring(0, $text.Length - 4) + "_FV2404"
        } elseif ($text.EndsWith("_new")) {
            $cell.Value2 = $text.Substring(0, $text.Length - 4) + "_FV2410"
        }
    }
}

# 2. Freeze the header row (pane split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the data range into an Excel Table ("Table1") with a header row
$rng = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
